$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Season" (N/O) columns ------------------------------------

# Header row (row 1): plain "Train" / "Test" labels, same as F1/G1 (no style)
$ws.Range("N1").Value = "Train"
$ws.Range("O1").Value = "Test"

# Row 2 (KNN): copy number-style from L2:M2, then overwrite values
$ws.Range("L2:M2").Copy()
$ws.Range("N2:O2").PasteSpecial(-4122)
$ws.Range("N2").Value = 0.76100000000000001
$ws.Range("O2").Value = 0.33700000000000002

# Row 3 (RF): copy number-style from L3:M3, then overwrite values
$ws.Range("L3:M3").Copy()
$ws.Range("N3:O3").PasteSpecial(-4122)
$ws.Range("N3").Value = 0.97099999999999997
$ws.Range("O3").Value = 0.48399999999999999

# Row 4 (SVR): copy number-style from L4:M4, then overwrite values
$ws.Range("L4:M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)
$ws.Range("N4").Value = 0.79600000000000004
$ws.Range("O4").Value = 0.496

$excel.CutCopyMode = 0

# --- New notes box for the "Total Generation" run (N7:O13) -------------

# N7 gets the bordered "box title" look (same as B7, i.e. style with
# borderId 4 / center / wrap) and the new note text.
$ws.Range("B7").Copy()
$ws.Range("N7:O7").PasteSpecial(-4122)
$ws.Range("N7").Value = "Total Generation (after bug fix) 0.85 - 0.15 split Same as previous"

# N8:O13 get the plain "box body" look (same as B8, i.e. center/wrap, no
# border), matching the other note boxes.
$ws.Range("B8").Copy()
$ws.Range("N8:O13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Merge the new note box, like the other five boxes in row 7:13
$ws.Range("N7:O13").Merge()

# --- Selection / view state ---------------------------------------------

$ws.Range("N5").Select()
